# Apply benefits-totals updates to the "DFT ID Number" worksheet.
# The sheet contains a grid of forecast cost rows; several cells move
# between "blank" and "0" and several numeric values are recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DFT ID Number")

# Row 23 - 18-19 RDEL Forecast one off new costs
$ws.Range("C23").Value = 0
$ws.Range("D23").ClearContents()
$ws.Range("E23").Value = 1.08

# Row 24 - 18-19 RDEL Forecast recurring new costs
$ws.Range("C24").Value = 0
$ws.Range("D24").ClearContents()
$ws.Range("E24").Value = 0

# Row 25 - 18-19 RDEL Forecast recurring old costs
$ws.Range("C25").Value = 0
$ws.Range("D25").ClearContents()

# Row 27 - 18-19 RDEL Forecast Total
$ws.Range("E27").Value = 1.08

# Row 29 - 18-19 CDEL Forecast one off new costs
$ws.Range("C29").Value = 447.8
$ws.Range("E29").Value = 14.25
$ws.Range("F29").Value = 78.40000000000001

# Row 30 - 18-19 CDEL Forecast recurring new costs
$ws.Range("D30").ClearContents()
$ws.Range("E30").Value = 0

# Row 31 - 18-19 CDEL Forecast recurring old costs
$ws.Range("D31").ClearContents()

# Row 32 - 18-19 Forecast Non-Gov
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 3.4
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 99.90000000000001

# Row 34 - 18-19 Forecast - Income both Revenue and Capital
$ws.Range("C34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0

# Row 35 - 19-20 RDEL Forecast one off new costs
$ws.Range("C35").Value = 0
$ws.Range("D35").ClearContents()
$ws.Range("E35").Value = 2

# Row 36 - 19-20 RDEL Forecast recurring new costs
$ws.Range("C36").Value = 0
$ws.Range("D36").ClearContents()
$ws.Range("E36").Value = 0

# Row 37 - 19-20 RDEL Forecast recurring old costs
$ws.Range("C37").Value = 0
$ws.Range("D37").ClearContents()

# Row 39 - 19-20 RDEL Forecast Total
$ws.Range("E39").Value = 2

# Row 41 - 19-20 CDEL Forecast one off new costs
$ws.Range("C41").Value = 320.9
$ws.Range("F41").Value = 165

# Row 42 - 19-20 CDEL Forecast recurring new costs
$ws.Range("D42").ClearContents()

# Row 43 - 19-20 CDEL Forecast recurring old costs
$ws.Range("D43").ClearContents()

# Row 44 - 19-20 Forecast Non-Gov
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 22.6
$ws.Range("E44").Value = 0
$ws.Range("G44").Value = 11.7

# Row 46 - 19-20 Forecast - Income both Revenue and Capital
$ws.Range("C46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
